$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.NumberFormat = "General"
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Cells.Item(2,4) "98.054.00"
Set-TextValue $ws.Cells.Item(2,5) "  -0.45%  "
# Row 3
Set-TextValue $ws.Cells.Item(3,4) "3.409.45"
Set-TextValue $ws.Cells.Item(3,5) "  +2.02%  "
# Row 4
Set-TextValue $ws.Cells.Item(4,5) "  -0.12%  "
# Row 5
Set-TextValue $ws.Cells.Item(5,4) "256.34"
Set-TextValue $ws.Cells.Item(5,5) "  -2.51%  "
# Row 6
Set-TextValue $ws.Cells.Item(6,4) "657.77"
Set-TextValue $ws.Cells.Item(6,5) "  +1.88%  "
# Row 7
Set-TextValue $ws.Cells.Item(7,4) "1.45"
Set-TextValue $ws.Cells.Item(7,5) "  -4.02%  "
# Row 8
Set-TextValue $ws.Cells.Item(8,4) "0.434"
Set-TextValue $ws.Cells.Item(8,5) "  -5.34%  "
# Row 9
Set-TextValue $ws.Cells.Item(9,4) "1.06"
Set-TextValue $ws.Cells.Item(9,5) "  -0.78%  "
# Row 10
Set-TextValue $ws.Cells.Item(10,5) "  -0.09%  "
# Row 11
Set-TextValue $ws.Cells.Item(11,4) "3.404.98"
Set-TextValue $ws.Cells.Item(11,5) "  +1.97%  "
# Row 12
Set-TextValue $ws.Cells.Item(12,4) "0.215"
Set-TextValue $ws.Cells.Item(12,5) "  +3.61%  "
# Row 13
Set-TextValue $ws.Cells.Item(13,4) "41.89"
Set-TextValue $ws.Cells.Item(13,5) "  -5.10%  "
# Row 14
Set-TextValue $ws.Cells.Item(14,4) "6.40"
Set-TextValue $ws.Cells.Item(14,5) "  +15.41%  "
# Row 15
Set-TextValue $ws.Cells.Item(15,4) "97.733.86"
Set-TextValue $ws.Cells.Item(15,5) "  -0.57%  "
# Row 16
Set-TextValue $ws.Cells.Item(16,4) "0.0000265"
Set-TextValue $ws.Cells.Item(16,5) "  -2.52%  "
# Row 17
Set-TextValue $ws.Cells.Item(17,4) "4.042.52"
Set-TextValue $ws.Cells.Item(17,5) "  +1.58%  "
# Row 18
Set-TextValue $ws.Cells.Item(18,4) "9.05"
Set-TextValue $ws.Cells.Item(18,5) "  +22.08%  "
# Row 19
Set-TextValue $ws.Cells.Item(19,4) "0.588"
Set-TextValue $ws.Cells.Item(19,5) "  +39.45%  "
# Row 20
Set-TextValue $ws.Cells.Item(20,4) "3.418.13"
Set-TextValue $ws.Cells.Item(20,5) "  +2.22%  "
# Row 21
Set-TextValue $ws.Cells.Item(21,4) "17.67"
Set-TextValue $ws.Cells.Item(21,5) "  +6.30%  "
# Row 22
Set-TextValue $ws.Cells.Item(22,4) "10.86"
Set-TextValue $ws.Cells.Item(22,5) "  +8.21%  "
# Row 23
$ws.Cells.Item(23,2).Value = "BitcoinCash"
$ws.Cells.Item(23,3).Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue $ws.Cells.Item(23,4) "514.31"
Set-TextValue $ws.Cells.Item(23,5) "  -3.08%  "
# Row 24
$ws.Cells.Item(24,2).Value = "SuiNetwork"
$ws.Cells.Item(24,3).Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-TextValue $ws.Cells.Item(24,4) "3.45"
Set-TextValue $ws.Cells.Item(24,5) "  -2.72%  "
# Row 25
Set-TextValue $ws.Cells.Item(25,4) "0.0000206"
Set-TextValue $ws.Cells.Item(25,5) "  -2.89%  "
# Row 26
Set-TextValue $ws.Cells.Item(26,4) "6.37"
Set-TextValue $ws.Cells.Item(26,5) "  +5.42%  "
# Row 27
Set-TextValue $ws.Cells.Item(27,4) "100.68"
Set-TextValue $ws.Cells.Item(27,5) "  -0.71%  "
# Row 28
Set-TextValue $ws.Cells.Item(28,4) "12.96"
Set-TextValue $ws.Cells.Item(28,5) "  +2.02%  "
# Row 29
Set-TextValue $ws.Cells.Item(29,4) "3.604.58"
Set-TextValue $ws.Cells.Item(29,5) "  +2.24%  "
# Row 30
Set-TextValue $ws.Cells.Item(30,4) "0.150"
Set-TextValue $ws.Cells.Item(30,5) "  +1.75%  "
# Row 31
Set-TextValue $ws.Cells.Item(31,4) "11.77"
Set-TextValue $ws.Cells.Item(31,5) "  +8.77%  "
# Row 32
Set-TextValue $ws.Cells.Item(32,4) "0.198"
Set-TextValue $ws.Cells.Item(32,5) "  +5.14%  "
# Row 33
Set-TextValue $ws.Cells.Item(33,4) "0.997"
Set-TextValue $ws.Cells.Item(33,5) "  -0.38%  "
# Row 34
Set-TextValue $ws.Cells.Item(34,4) "0.578"
Set-TextValue $ws.Cells.Item(34,5) "  +12.91%  "
# Row 35
$ws.Cells.Item(35,2).Value = "Binance-PegBSC-USD"
$ws.Cells.Item(35,3).Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue $ws.Cells.Item(35,4) "1.00"
Set-TextValue $ws.Cells.Item(35,5) "  -0.97%  "
# Row 36
$ws.Cells.Item(36,2).Value = "PancakeSwap"
$ws.Cells.Item(36,3).Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws.Cells.Item(36,4) "2.36"
Set-TextValue $ws.Cells.Item(36,5) "  +15.09%  "
# Row 37
Set-TextValue $ws.Cells.Item(37,4) "29.93"
Set-TextValue $ws.Cells.Item(37,5) "  +3.10%  "
# Row 38
Set-TextValue $ws.Cells.Item(38,4) "7.84"
Set-TextValue $ws.Cells.Item(38,5) "  +0.95%  "
# Row 39
Set-TextValue $ws.Cells.Item(39,4) "1.46"
Set-TextValue $ws.Cells.Item(39,5) "  +10.91%  "
# Row 40
Set-TextValue $ws.Cells.Item(40,4) "530.83"
Set-TextValue $ws.Cells.Item(40,5) "  +1.90%  "
# Row 41
Set-TextValue $ws.Cells.Item(41,4) "0.152"
Set-TextValue $ws.Cells.Item(41,5) "  -1.99%  "
# Row 42
Set-TextValue $ws.Cells.Item(42,5) "  -0.01%  "
# Row 43
Set-TextValue $ws.Cells.Item(43,4) "0.877"
Set-TextValue $ws.Cells.Item(43,5) "  +8.94%  "
# Row 44
$ws.Cells.Item(44,2).Value = "WhiteBITCoin"
$ws.Cells.Item(44,3).Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextValue $ws.Cells.Item(44,4) "24.70"
Set-TextValue $ws.Cells.Item(44,5) "  -0.04%  "
# Row 45
$ws.Cells.Item(45,2).Value = "Cosmos"
$ws.Cells.Item(45,3).Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws.Cells.Item(45,4) "9.09"
Set-TextValue $ws.Cells.Item(45,5) "  +18.50%  "
# Row 46
$ws.Cells.Item(46,2).Value = "Filecoin"
$ws.Cells.Item(46,3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Cells.Item(46,4) "5.88"
Set-TextValue $ws.Cells.Item(46,5) "  +20.11%  "
# Row 47
Set-TextValue $ws.Cells.Item(47,4) "0.0428"
Set-TextValue $ws.Cells.Item(47,5) "  +10.50%  "
# Row 48
$ws.Cells.Item(48,2).Value = "MantraDAO"
$ws.Cells.Item(48,3).Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
Set-TextValue $ws.Cells.Item(48,4) "3.73"
Set-TextValue $ws.Cells.Item(48,5) "  -3.69%  "
# Row 49
$ws.Cells.Item(49,2).Value = "ImmutableX"
$ws.Cells.Item(49,3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Cells.Item(49,4) "1.68"
Set-TextValue $ws.Cells.Item(49,5) "  +15.58%  "
# Row 50
Set-TextValue $ws.Cells.Item(50,4) "3.31"
Set-TextValue $ws.Cells.Item(50,5) "  -0.55%  "
# Row 51
Set-TextValue $ws.Cells.Item(51,4) "2.11"
Set-TextValue $ws.Cells.Item(51,5) "  +4.88%  "
